$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at 449-452, pushing existing rows 449-481 down to 453-485.
$ws.Range("A449:R452").EntireRow.Insert()

# Populate the 4 newly inserted rows with the new weekly price entries
# (same Mercado/Region/Categoria/Variedad/Unidad/Origen/Kg as all other rows).

# Row 449 - Especial
$ws.Range("A449").Value = 2
$ws.Range("B449").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C449").Value = "Coquimbo"
$ws.Range("D449").Value = 44783
$ws.Range("E449").Value = 4
$ws.Range("F449").Value = 100112043
$ws.Range("G449").Value = "Pepino dulce"
$ws.Range("H449").Value = "Cultivar IV Región"
$ws.Range("I449").Value = "Especial"
$ws.Range("J449").Value = 248
$ws.Range("K449").Value = 12000
$ws.Range("L449").Value = 13000
$ws.Range("M449").Value = 12484
$ws.Range("N449").Value = "`$/bandeja 18 kilos"
$ws.Range("O449").Value = "Provincia de Limarí"
$ws.Range("P449").Value = 694
$ws.Range("Q449").Value = 18
$ws.Range("R449").Value = "Hortaliza"

# Row 450 - Primera
$ws.Range("A450").Value = 2
$ws.Range("B450").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C450").Value = "Coquimbo"
$ws.Range("D450").Value = 44783
$ws.Range("E450").Value = 4
$ws.Range("F450").Value = 100112043
$ws.Range("G450").Value = "Pepino dulce"
$ws.Range("H450").Value = "Cultivar IV Región"
$ws.Range("I450").Value = "Primera"
$ws.Range("J450").Value = 500
$ws.Range("K450").Value = 10000
$ws.Range("L450").Value = 11000
$ws.Range("M450").Value = 10500
$ws.Range("N450").Value = "`$/bandeja 18 kilos"
$ws.Range("O450").Value = "Provincia de Limarí"
$ws.Range("P450").Value = 583
$ws.Range("Q450").Value = 18
$ws.Range("R450").Value = "Hortaliza"

# Row 451 - Segunda
$ws.Range("A451").Value = 2
$ws.Range("B451").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C451").Value = "Coquimbo"
$ws.Range("D451").Value = 44783
$ws.Range("E451").Value = 4
$ws.Range("F451").Value = 100112043
$ws.Range("G451").Value = "Pepino dulce"
$ws.Range("H451").Value = "Cultivar IV Región"
$ws.Range("I451").Value = "Segunda"
$ws.Range("J451").Value = 300
$ws.Range("K451").Value = 8000
$ws.Range("L451").Value = 9000
$ws.Range("M451").Value = 8500
$ws.Range("N451").Value = "`$/bandeja 18 kilos"
$ws.Range("O451").Value = "Provincia de Limarí"
$ws.Range("P451").Value = 472
$ws.Range("Q451").Value = 18
$ws.Range("R451").Value = "Hortaliza"

# Row 452 - Tercera
$ws.Range("A452").Value = 2
$ws.Range("B452").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C452").Value = "Coquimbo"
$ws.Range("D452").Value = 44783
$ws.Range("E452").Value = 4
$ws.Range("F452").Value = 100112043
$ws.Range("G452").Value = "Pepino dulce"
$ws.Range("H452").Value = "Cultivar IV Región"
$ws.Range("I452").Value = "Tercera"
$ws.Range("J452").Value = 240
$ws.Range("K452").Value = 5000
$ws.Range("L452").Value = 6000
$ws.Range("M452").Value = 5500
$ws.Range("N452").Value = "`$/bandeja 18 kilos"
$ws.Range("O452").Value = "Provincia de Limarí"
$ws.Range("P452").Value = 306
$ws.Range("Q452").Value = 18
$ws.Range("R452").Value = "Hortaliza"
